$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'88.199.62"
$ws.Range("E2").Value = "'  -0.87%  "
$ws.Range("D3").Value = "'3.036.82"
$ws.Range("E3").Value = "'  -1.74%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'208.65"
$ws.Range("E5").Value = "'  -1.61%  "
$ws.Range("D6").Value = "'606.45"
$ws.Range("E6").Value = "'  -3.16%  "
$ws.Range("E7").Value = "'  -7.26%  "
$ws.Range("D8").Value = "'0.874"
$ws.Range("E8").Value = "'  +22.90%  "
$ws.Range("E9").Value = "'  +0.06%  "
$ws.Range("D10").Value = "'3.036.28"
$ws.Range("E10").Value = "'  -1.63%  "
$ws.Range("D11").Value = "'0.645"
$ws.Range("E11").Value = "'  +17.86%  "
$ws.Range("E12").Value = "'  +4.16%  "
$ws.Range("E13").Value = "'  -5.11%  "
$ws.Range("D14").Value = "'5.35"
$ws.Range("E14").Value = "'  +2.65%  "
$ws.Range("D15").Value = "'88.276.46"
$ws.Range("E15").Value = "'  -0.43%  "
$ws.Range("D16").Value = "'3.600.24"
$ws.Range("E16").Value = "'  -1.73%  "
$ws.Range("D17").Value = "'31.64"
$ws.Range("E17").Value = "'  -0.11%  "
$ws.Range("D18").Value = "'3.068.40"
$ws.Range("E18").Value = "'  -0.66%  "
$ws.Range("D19").Value = "'3.34"
$ws.Range("E19").Value = "'  +1.28%  "
$ws.Range("D20").Value = "'0.0000204"
$ws.Range("E20").Value = "'  -0.47%  "
$ws.Range("D21").Value = "'13.28"
$ws.Range("E21").Value = "'  +2.47%  "
$ws.Range("D22").Value = "'419.66"
$ws.Range("E22").Value = "'  -0.42%  "
$ws.Range("E23").Value = "'  +2.13%  "
$ws.Range("D24").Value = "'8.02"
$ws.Range("E24").Value = "'  -2.18%  "
$ws.Range("E25").Value = "'  +3.10%  "
$ws.Range("D26").Value = "'83.13"
$ws.Range("E26").Value = "'  +6.15%  "
$ws.Range("D27").Value = "'11.46"
$ws.Range("E27").Value = "'  +2.66%  "
$ws.Range("D28").Value = "'3.205.14"
$ws.Range("E28").Value = "'  -2.38%  "
$ws.Range("E29").Value = "'  -0.15%  "
$ws.Range("D30").Value = "'1.09"
$ws.Range("E30").Value = "'  +9.06%  "
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "'  +2.64%  "
$ws.Range("D32").Value = "'8.13"
$ws.Range("E32").Value = "'  +0.55%  "
$ws.Range("D33").Value = "'498.39"
$ws.Range("E33").Value = "'  -0.22%  "
$ws.Range("E34").Value = "'  -7.36%  "
$ws.Range("E35").Value = "'  -2.48%  "
$ws.Range("E36").Value = "'  -1.51%  "
$ws.Range("D37").Value = "'22.37"
$ws.Range("E37").Value = "'  +3.60%  "
$ws.Range("E38").Value = "'  -1.53%  "
$ws.Range("D39").Value = "'22.17"
$ws.Range("E39").Value = "'  -0.11%  "
$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "'  +3.97%  "
$ws.Range("E41").Value = "'  +0.24%  "
$ws.Range("E42").Value = "'  -0.04%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "'  +11.67%  "
$ws.Range("D44").Value = "'0.360"
$ws.Range("E44").Value = "'  +0.55%  "
$ws.Range("D45").Value = "'146.59"
$ws.Range("E45").Value = "'  +1.63%  "
$ws.Range("D46").Value = "'1.80"
$ws.Range("E46").Value = "'  -2.41%  "
$ws.Range("D47").Value = "'43.34"
$ws.Range("E47").Value = "'  +0.19%  "
$ws.Range("D48").Value = "'0.0684"
$ws.Range("E48").Value = "'  +14.28%  "
$ws.Range("E49").Value = "'  +3.18%  "
$ws.Range("E50").Value = "'  +2.67%  "
$ws.Range("D51").Value = "'154.38"
$ws.Range("E51").Value = "'  -5.90%  "
